{"js": "// Replace the AI-automation abstract paragraph with the new Power-Plant /\n// linear-regression abstract paragraph, keeping the paragraph's existing\n// formatting (black text color, no-border/no-first-line-indent pPr, etc.)\n// intact \u2014 only the paragraph's textual content changes.\n\nconst OLD_SNIPPET = \"Artificial Intelligence (AI) is one of the most transformative\";\n\nconst NEW_TEXT =\n  \"Since the Twentieth Century, Electric Power has been the source of advancement of Human Technology. \" +\n  \"It has been a basic need for the past 100 years. The Power Plants of today are significantly more efficient than ever, \" +\n  \"like the Combined Cycle Power Plants which take advantage of the heat exhaust \" +\n  \"generated by their predecessor designs and reuse it to generate even more electricity. The level of efficiency is so high that we can now measure the effects of the environment on the productivity of the plant. \" +\n  \"Measuring the variance in the output of an electric plant is essential not only for the workers of the facility, but to the investors, the end users, and to the Economy. \" +\n  \"This paper explores how linear regression can be used to predict the output of a power plant based off of ambient temperature, atmospheric pressure, and other environment variables \" +\n  \"automatically with a prediction model.\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.indexOf(OLD_SNIPPET) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the abstract paragraph to replace.\");\n}\n\n// \"Replace\" swaps the paragraph's full text while preserving the run\n// formatting (w:rPr -> color 000000) already on it.\ntarget.insertText(NEW_TEXT, \"Replace\");\nawait context.sync();\n", "ps1": "# Replace the AI-automation abstract paragraph with the new Power-Plant /\n# linear-regression abstract paragraph, keeping the paragraph's existing\n# formatting (black text color, no-border/no-first-line-indent pPr, etc.)\n# intact -- only the paragraph's textual content changes.\n\n$oldSnippet = \"Artificial Intelligence (AI) is one of the most transformative\"\n\n$newText = \"Since the Twentieth Century, Electric Power has been the source of advancement of Human Technology. \" + `\n  \"It has been a basic need for the past 100 years. The Power Plants of today are significantly more efficient than ever, \" + `\n  \"like the Combined Cycle Power Plants which take advantage of the heat exhaust \" + `\n  \"generated by their predecessor designs and reuse it to generate even more electricity. The level of efficiency is so high that we can now measure the effects of the environment on the productivity of the plant. \" + `\n  \"Measuring the variance in the output of an electric plant is essential not only for the workers of the facility, but to the investors, the end users, and to the Economy. \" + `\n  \"This paper explores how linear regression can be used to predict the output of a power plant based off of ambient temperature, atmospheric pressure, and other environment variables \" + `\n  \"automatically with a prediction model.\"\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($oldSnippet)) {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the abstract paragraph to replace.\"\n}\n\n# Replace only the paragraph's text (not its trailing paragraph mark), so the\n# paragraph formatting (w:pPr) and run formatting (w:rPr -> color 000000)\n# already on the run are preserved.\n$paraRange = $target.Range\n$textRange = $d.Range($paraRange.Start, $paraRange.End - 1)\n$textRange.Text = $newText\n"}
